$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; existing rows 10-35 shift down to 11-36.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 45014
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112044
$ws.Range("G10").Value = "Perejil"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 1500
$ws.Range("L10").Value = 1500
$ws.Range("M10").Value = 1500
$ws.Range("N10").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O10").Value = "Región del Maule"
$ws.Range("P10").Value = 1500
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"

# Make sure the date cell keeps the workbook's date number format (same as other D cells).
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
